# BAU Ban New Power Plants This Year.xlsx
# Commit: Update historical capacity additions; for BPMCCS, increase wind and
# solar additions to cover all capacity currently under construction.
#
# Concretely, on the "BBNPPTY" sheet, 2024 (column E) is now treated the same
# as a historical year (2021-2023, columns B-D) for every technology: the
# "allowed" flag is set to 1 (instead of mostly 0) and the cell's number
# formatting is cleared so it matches the unformatted look of the other
# historical-year columns. A new explanatory note is added to the "About"
# sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BBNPPTY")

# --- BBNPPTY sheet: set 2024 (column E) allowed-flag to 1 for every row ---
# Rows 2-14 and 16-25 hold technology flags; row 1 is the year header and
# row 15 (offshore wind) is already 1 with no special formatting.
$rngTop = $ws2.Range("E2:E14")
$rngTop.ClearFormats()
$rngTop.Value = 1

$rngBottom = $ws2.Range("E16:E25")
$rngBottom.ClearFormats()
$rngBottom.Value = 1

# --- About sheet: add a new note explaining the historical-year treatment ---
$ws1.Range("A22").Value = "For historical years (2021-2024), we directly read in capacity additions in other input"
$ws1.Range("A23").Value = "data files and therefore disallow additional economic additions."

# --- Restore each sheet's own selection state (About stays the active tab) ---
$ws2.Range("D2:E25").Select()
$ws1.Range("A24").Select()
